$wb = $excel.ActiveWorkbook

# This script applies updated market-price figures (columns H-N) produced by
# the scheduled Kraken market-data refresh, row by row, per sheet. Where a
# cell has no new value (no market data available for that sub-column), it is
# cleared to $null so the cell is omitted entirely, matching the refreshed feed.

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")

# Row 2
$ws.Range("H2").Value = 342.44446
$ws.Range("I2").Value = 335.25
$ws.Range("K2").Value = 335.25
$ws.Range("M2").Value = -222.25

# Row 17
$ws.Range("H17").Value = 1481.1333
$ws.Range("I17").Value = 1249.75
$ws.Range("J17").Value = 1565.2727
$ws.Range("K17").Value = 3749.25
$ws.Range("L17").Value = 4695.8181
$ws.Range("M17").Value = -3581.25
$ws.Range("N17").Value = -5031.8181

# Row 33
$ws.Range("H33").Value = 192.33333
$ws.Range("I33").Value = 192.33333
$ws.Range("K33").Value = 192.33333
$ws.Range("M33").Value = 36.66667000000001

# Row 40
$ws.Range("H40").Value = 4500
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 4500
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 4500
$ws.Range("M40").Value = $null
$ws.Range("N40").Value = -4850

# Row 69
$ws.Range("H69").Value = 0
$ws.Range("I69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("M69").Value = $null

# Row 72
$ws.Range("H72").Value = 0
$ws.Range("I72").Value = 0
$ws.Range("K72").Value = 0
$ws.Range("M72").Value = $null

# Row 121
$ws.Range("H121").Value = 1999.5
$ws.Range("J121").Value = 1999.5
$ws.Range("L121").Value = 5998.5
$ws.Range("N121").Value = -9492.5

# Row 135
$ws.Range("H135").Value = 1497.5555
$ws.Range("I135").Value = 1068.4286
$ws.Range("K135").Value = 9615.857399999999
$ws.Range("M135").Value = -7080.857399999999

# Row 137
$ws.Range("H137").Value = 1966.6666
$ws.Range("I137").Value = 2500
$ws.Range("J137").Value = 900
$ws.Range("K137").Value = 7500
$ws.Range("L137").Value = 2700
$ws.Range("M137").Value = -4950
$ws.Range("N137").Value = -7800

# Row 138
$ws.Range("H138").Value = 3360.375
$ws.Range("I138").Value = 2697.5454
$ws.Range("K138").Value = 8092.6362
$ws.Range("M138").Value = -2952.6362

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")

# Row 32
$ws.Range("H32").Value = 3267
$ws.Range("I32").Value = 3311.0454
$ws.Range("K32").Value = 3311.0454
$ws.Range("M32").Value = -3024.0454

# Row 122
$ws.Range("H122").Value = 4066.4119
$ws.Range("I122").Value = 4190.6924
$ws.Range("J122").Value = 3662.5
$ws.Range("K122").Value = 12572.0772
$ws.Range("L122").Value = 10987.5
$ws.Range("M122").Value = -10122.0772
$ws.Range("N122").Value = -15887.5

# Row 132
$ws.Range("H132").Value = 4896.6
$ws.Range("I132").Value = 3911.1667
$ws.Range("J132").Value = 6374.75
$ws.Range("K132").Value = 11733.5001
$ws.Range("L132").Value = 19124.25
$ws.Range("M132").Value = -9203.500100000001
$ws.Range("N132").Value = -24184.25

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")

# Row 11
$ws.Range("H11").Value = 1003.8333
$ws.Range("I11").Value = 1268.25
$ws.Range("J11").Value = 475
$ws.Range("K11").Value = 1268.25
$ws.Range("L11").Value = 475
$ws.Range("M11").Value = -1128.25
$ws.Range("N11").Value = -755

# Row 19
$ws.Range("H19").Value = 110
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 110
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 110
$ws.Range("M19").Value = $null
$ws.Range("N19").Value = -456

# Row 20
$ws.Range("H20").Value = 2179.7144
$ws.Range("I20").Value = 1702
$ws.Range("J20").Value = 2816.6667
$ws.Range("K20").Value = 1702
$ws.Range("L20").Value = 2816.6667
$ws.Range("M20").Value = -1455
$ws.Range("N20").Value = -3310.6667

# Row 76
$ws.Range("H76").Value = 7657
$ws.Range("J76").Value = 7657
$ws.Range("L76").Value = 7657
$ws.Range("N76").Value = -8287

# Row 79
$ws.Range("H79").Value = 7657
$ws.Range("J79").Value = 7657
$ws.Range("L79").Value = 7657
$ws.Range("N79").Value = -9841

# Row 86
$ws.Range("H86").Value = 1355
$ws.Range("I86").Value = 1261.6666
$ws.Range("J86").Value = 1495
$ws.Range("K86").Value = 1261.6666
$ws.Range("L86").Value = 1495
$ws.Range("M86").Value = -138.6666
$ws.Range("N86").Value = -3741

# Row 89
$ws.Range("H89").Value = 1355
$ws.Range("I89").Value = 1261.6666
$ws.Range("J89").Value = 1495
$ws.Range("K89").Value = 6308.333000000001
$ws.Range("L89").Value = 7475
$ws.Range("M89").Value = -692.3330000000005
$ws.Range("N89").Value = -18707

# Row 105
$ws.Range("H105").Value = 4059.6
$ws.Range("I105").Value = 4749.5
$ws.Range("J105").Value = 3599.6667
$ws.Range("K105").Value = 4749.5
$ws.Range("L105").Value = 3599.6667
$ws.Range("M105").Value = -3002.5
$ws.Range("N105").Value = -7093.6667

# Row 134
$ws.Range("H134").Value = 6599.1665
$ws.Range("I134").Value = 4399.5
$ws.Range("J134").Value = 10998.5
$ws.Range("K134").Value = 13198.5
$ws.Range("L134").Value = 32995.5
$ws.Range("M134").Value = -10663.5
$ws.Range("N134").Value = -38065.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")

# Row 7
$ws.Range("H7").Value = 665
$ws.Range("J7").Value = 665
$ws.Range("L7").Value = 665
$ws.Range("N7").Value = -891

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")

# Row 4
$ws.Range("H4").Value = 546182.8
$ws.Range("I4").Value = 500801.2
$ws.Range("K4").Value = 1502403.6
$ws.Range("M4").Value = -1502291.6

# Row 68
$ws.Range("H68").Value = 1051.3
$ws.Range("I68").Value = 1283.3334
$ws.Range("K68").Value = 3850.0002
$ws.Range("M68").Value = -3039.0002

# Row 71
$ws.Range("H71").Value = 1051.3
$ws.Range("I71").Value = 1283.3334
$ws.Range("K71").Value = 11550.0006
$ws.Range("M71").Value = -7494.000599999999

# Row 92
$ws.Range("H92").Value = 5000
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").Value = $null

# Row 97
$ws.Range("H97").Value = 700
$ws.Range("I97").Value = 700
$ws.Range("K97").Value = 2100
$ws.Range("M97").Value = -1604

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")

# Row 7
$ws.Range("H7").Value = 254100.75
$ws.Range("I7").Value = 336134.34
$ws.Range("J7").Value = 8000
$ws.Range("K7").Value = 336134.34
$ws.Range("L7").Value = 8000
$ws.Range("M7").Value = -336022.34
$ws.Range("N7").Value = -8224

# Row 8
$ws.Range("H8").Value = 254100.75
$ws.Range("I8").Value = 336134.34
$ws.Range("J8").Value = 8000
$ws.Range("K8").Value = 336134.34
$ws.Range("L8").Value = 8000
$ws.Range("M8").Value = -335995.34
$ws.Range("N8").Value = -8278

# Row 14
$ws.Range("H14").Value = 14917289
$ws.Range("I14").Value = 20708916
$ws.Range("J14").Value = 3334035
$ws.Range("K14").Value = 20708916
$ws.Range("L14").Value = 3334035
$ws.Range("M14").Value = -20708748
$ws.Range("N14").Value = -3334371

# Row 21
$ws.Range("H21").Value = 2000000
$ws.Range("I21").Value = 2000000
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 2000000
$ws.Range("L21").Value = 0
$ws.Range("M21").Value = -1999827
$ws.Range("N21").Value = $null

# Row 30
$ws.Range("H30").Value = 2000000
$ws.Range("I30").Value = 2000000
$ws.Range("J30").Value = 0
$ws.Range("K30").Value = 2000000
$ws.Range("L30").Value = 0
$ws.Range("M30").Value = -1999895
$ws.Range("N30").Value = $null

# Row 70
$ws.Range("H70").Value = 3630
$ws.Range("I70").Value = 2945
$ws.Range("J70").Value = 5000
$ws.Range("K70").Value = 2945
$ws.Range("L70").Value = 5000
$ws.Range("M70").Value = -2675
$ws.Range("N70").Value = -5540

# Row 73
$ws.Range("H73").Value = 3630
$ws.Range("I73").Value = 2945
$ws.Range("J73").Value = 5000
$ws.Range("K73").Value = 2945
$ws.Range("L73").Value = 5000
$ws.Range("M73").Value = -2009
$ws.Range("N73").Value = -6872

# Row 102
$ws.Range("H102").Value = 1547.5714
$ws.Range("I102").Value = 1547.5714
$ws.Range("K102").Value = 1547.5714
$ws.Range("M102").Value = 74.42859999999996

# Row 134
$ws.Range("H134").Value = 99750
$ws.Range("J134").Value = 99750
$ws.Range("L134").Value = 299250
$ws.Range("N134").Value = -304320

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")

# Row 132
$ws.Range("H132").Value = 5500
$ws.Range("I132").Value = 5500
$ws.Range("K132").Value = 16500
$ws.Range("M132").Value = -13970

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")

# Row 132
$ws.Range("H132").Value = 7034.385
$ws.Range("I132").Value = 5888.6665
$ws.Range("J132").Value = 9612.25
$ws.Range("K132").Value = 17665.9995
$ws.Range("L132").Value = 28836.75
$ws.Range("M132").Value = -15135.9995
$ws.Range("N132").Value = -33896.75

# Row 140
$ws.Range("H140").Value = 80000
$ws.Range("J140").Value = 80000
$ws.Range("L140").Value = 80000
$ws.Range("N140").Value = -90360
